$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 237 ("Ajo - Chino - Primera",
# fecha 2021-12-13 / serial 44543). All the existing rows from 237 downward
# shift down by one (old row 296 becomes new row 297).
$ws.Rows.Item(237).Insert()

$ws.Range("A237").Value = 3
$ws.Range("B237").Value = "Femacal de La Calera"
$ws.Range("C237").Value = "Coquimbo"
$ws.Range("D237").Value = 44543
$ws.Range("E237").Value = 5
$ws.Range("F237").Value = 100112003
$ws.Range("G237").Value = "Ajo"
$ws.Range("H237").Value = "Chino"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 73
$ws.Range("K237").Value = 15500
$ws.Range("L237").Value = 16000
$ws.Range("M237").Value = 15740
$ws.Range("N237").Value = "$/caja 10 kilos"
$ws.Range("O237").Value = "China"
$ws.Range("P237").Value = 1574
$ws.Range("Q237").Value = 10
$ws.Range("R237").Value = "Hortaliza"
